$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.694.19"
$ws.Range("E2").Value = "  +6.07%  "

$ws.Range("D3").Value = "2.312.47"
$ws.Range("E3").Value = "  +5.19%  "

$ws.Range("E4").Value = "  -0.62%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.69"
$ws.Range("E5").Value = "  +1.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.47"
$ws.Range("E6").Value = "  +10.62%  "

$ws.Range("E7").Value = "  -0.64%  "

$ws.Range("E8").Value = "  -0.55%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("E9").Value = "  +10.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.53"
$ws.Range("E10").Value = "  +9.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0799"
$ws.Range("E11").Value = "  +3.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.36"
$ws.Range("E12").Value = "  +9.37%  "

$ws.Range("E13").Value = "  +1.31%  "

$ws.Range("D14").Value = "2.663.40"
$ws.Range("E14").Value = "  +4.98%  "

$ws.Range("D15").Value = "2.313.94"
$ws.Range("E15").Value = "  +1.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.96"
$ws.Range("E16").Value = "  +6.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.820"
$ws.Range("E17").Value = "  +7.52%  "

$ws.Range("D18").Value = "46.636.70"
$ws.Range("E18").Value = "  +6.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.17"
$ws.Range("E19").Value = "  +22.47%  "

$ws.Range("D20").Value = "0.0₃0940"
$ws.Range("E20").Value = "  +6.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.15"
$ws.Range("E21").Value = "  +5.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.88"
$ws.Range("E22").Value = "  +6.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.36"
$ws.Range("E23").Value = "  +8.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.93"
$ws.Range("E24").Value = "  +6.10%  "

$ws.Range("E25").Value = "  +9.21%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.92"
$ws.Range("E27").Value = "  +20.68%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.24"
$ws.Range("E28").Value = "  +1.46%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.86"
$ws.Range("E29").Value = "  +7.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.11"
$ws.Range("E30").Value = "  +6.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.79"
$ws.Range("E31").Value = "  +9.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "147.86"
$ws.Range("E32").Value = "  +0.06%  "

$ws.Range("E33").Value = "  +9.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.63"
$ws.Range("E34").Value = "  +5.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.14"
$ws.Range("E35").Value = "  +10.60%  "

$ws.Range("E36").Value = "  +10.56%  "

$ws.Range("E37").Value = "  +2.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.81"
$ws.Range("E38").Value = "  +10.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.37"
$ws.Range("E39").Value = "  +15.50%  "

$ws.Range("E40").Value = "  +14.85%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.46"
$ws.Range("E41").Value = "  +13.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0308"
$ws.Range("E42").Value = "  +9.98%  "

$ws.Range("E43").Value = "  -0.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.00"
$ws.Range("E44").Value = "  +22.75%  "

$ws.Range("D45").Value = "1.842.24"
$ws.Range("E45").Value = "  +5.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.58"
$ws.Range("E46").Value = "  +23.67%  "

$ws.Range("E47").Value = "  +17.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "74.45"
$ws.Range("E48").Value = "  +10.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.97"
$ws.Range("E49").Value = "  +12.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.59"
$ws.Range("E50").Value = "  +6.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.22"
$ws.Range("E51").Value = "  +10.40%  "
